$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - fill in Item/Details first (Details before Item, matching authoring order),
# then Recorded date (copy date format from row above to avoid creating a new style),
# then Status.
$ws.Range("D6").Value = "Camera not showing after hosting to server and opening in client laptop"
$ws.Range("C6").Value = "Camera not displaying"
$ws.Range("E5").Copy($ws.Range("E6"))
$ws.Range("E6").Value2 = 44016
$ws.Range("G6").Value = "Pending"

# Row 7 - Item before Details.
$ws.Range("C7").Value = "On page refresh error"
$ws.Range("D7").Value = "Not found error showing if we refresh the page"
$ws.Range("E5").Copy($ws.Range("E7"))
$ws.Range("E7").Value2 = 44016
$ws.Range("G7").Value = "Pending"

# Row 8 - Item before Details.
$ws.Range("C8").Value = "Slow page"
$ws.Range("D8").Value = "View Service entry page is very slow"
$ws.Range("E5").Copy($ws.Range("E8"))
$ws.Range("E8").Value2 = 44016
$ws.Range("G8").Value = "Pending"

# Row heights for the newly wrapped rows.
$ws.Rows.Item(6).RowHeight = 29
$ws.Rows.Item(7).RowHeight = 29

# Update the active selection as recorded in the sheet view.
$ws.Range("F10").Select()
